$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "29.339.75"
$ws.Cells.Item(3,4).Value = "1.847.13"
$ws.Cells.Item(3,5).Value = "  -0.25%  "
$c = $ws.Cells.Item(4,4)
$c.NumberFormat = "@"
$c.Value = "0.9973"
$c.Style = "Normal"
$ws.Cells.Item(4,5).Value = "  -0.25%  "
$c = $ws.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value = "240.04"
$c.Style = "Normal"
$ws.Cells.Item(5,5).Value = "  -0.41%  "
$c = $ws.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value = "0.6264"
$c.Style = "Normal"
$ws.Cells.Item(6,5).Value = "  -0.54%  "
$c = $ws.Cells.Item(7,4)
$c.NumberFormat = "@"
$c.Value = "0.9987"
$c.Style = "Normal"
$ws.Cells.Item(7,5).Value = "  -0.18%  "
$c = $ws.Cells.Item(8,4)
$c.NumberFormat = "@"
$c.Value = "0.07596"
$c.Style = "Normal"
$ws.Cells.Item(8,5).Value = "  -1.18%  "
$c = $ws.Cells.Item(9,4)
$c.NumberFormat = "@"
$c.Value = "0.2907"
$c.Style = "Normal"
$ws.Cells.Item(9,5).Value = "  -1.09%  "
$c = $ws.Cells.Item(10,4)
$c.NumberFormat = "@"
$c.Value = "24.64"
$c.Style = "Normal"
$ws.Cells.Item(10,5).Value = "  +0.28%  "
$c = $ws.Cells.Item(11,4)
$c.NumberFormat = "@"
$c.Value = "0.07738"
$c.Style = "Normal"
$c = $ws.Cells.Item(12,4)
$c.NumberFormat = "@"
$c.Value = "5.019"
$c.Style = "Normal"
$ws.Cells.Item(12,5).Value = "  -0.24%  "
$c = $ws.Cells.Item(13,4)
$c.NumberFormat = "@"
$c.Value = "0.6786"
$c.Style = "Normal"
$ws.Cells.Item(13,5).Value = "  -0.41%  "
$ws.Cells.Item(14,5).Value = "  -3.59%  "
$c = $ws.Cells.Item(15,4)
$c.NumberFormat = "@"
$c.Value = "82.97"
$c.Style = "Normal"
$ws.Cells.Item(15,5).Value = "  -0.89%  "
$c = $ws.Cells.Item(16,4)
$c.NumberFormat = "@"
$c.Value = "6.126"
$c.Style = "Normal"
$ws.Cells.Item(16,5).Value = "  -0.41%  "
$ws.Cells.Item(17,4).Value = "29.349.59"
$ws.Cells.Item(17,5).Value = "  -0.36%  "
$c = $ws.Cells.Item(18,4)
$c.NumberFormat = "@"
$c.Value = "229.08"
$c.Style = "Normal"
$ws.Cells.Item(18,5).Value = "  -0.11%  "
$c = $ws.Cells.Item(19,4)
$c.NumberFormat = "@"
$c.Value = "12.32"
$c.Style = "Normal"
$ws.Cells.Item(19,5).Value = "  -1.31%  "
$c = $ws.Cells.Item(20,4)
$c.NumberFormat = "@"
$c.Value = "0.9986"
$c.Style = "Normal"
$ws.Cells.Item(20,5).Value = "  -0.22%  "
$c = $ws.Cells.Item(21,4)
$c.NumberFormat = "@"
$c.Value = "7.469"
$c.Style = "Normal"
$ws.Cells.Item(21,5).Value = "  +0.12%  "
$c = $ws.Cells.Item(22,4)
$c.NumberFormat = "@"
$c.Value = "0.9987"
$c.Style = "Normal"
$ws.Cells.Item(22,5).Value = "  -0.18%  "
$c = $ws.Cells.Item(23,4)
$c.NumberFormat = "@"
$c.Value = "158.45"
$c.Style = "Normal"
$ws.Cells.Item(23,5).Value = "  +0.82%  "
$ws.Cells.Item(24,5).Value = "  -0.33%  "
$ws.Cells.Item(25,5).Value = "  +0.47%  "
$c = $ws.Cells.Item(26,4)
$c.NumberFormat = "@"
$c.Value = "17.67"
$c.Style = "Normal"
$ws.Cells.Item(26,5).Value = "  -0.04%  "
$c = $ws.Cells.Item(27,4)
$c.NumberFormat = "@"
$c.Value = "1.434"
$c.Style = "Normal"
$ws.Cells.Item(27,5).Value = "  +9.24%  "
$ws.Cells.Item(28,5).Value = "  +0.06%  "
$c = $ws.Cells.Item(29,4)
$c.NumberFormat = "@"
$c.Value = "0.05600"
$c.Style = "Normal"
$ws.Cells.Item(29,5).Value = "  -2.06%  "
$c = $ws.Cells.Item(30,4)
$c.NumberFormat = "@"
$c.Value = "4.102"
$c.Style = "Normal"
$ws.Cells.Item(30,5).Value = "  -0.72%  "
$c = $ws.Cells.Item(31,4)
$c.NumberFormat = "@"
$c.Value = "4.063"
$c.Style = "Normal"
$ws.Cells.Item(31,5).Value = "  +0.26%  "
$c = $ws.Cells.Item(32,4)
$c.NumberFormat = "@"
$c.Value = "1.160"
$c.Style = "Normal"
$ws.Cells.Item(32,5).Value = "  -0.30%  "
$c = $ws.Cells.Item(33,4)
$c.NumberFormat = "@"
$c.Value = "1.829"
$c.Style = "Normal"
$ws.Cells.Item(33,5).Value = "  -1.27%  "
$c = $ws.Cells.Item(34,4)
$c.NumberFormat = "@"
$c.Value = "0.6954"
$c.Style = "Normal"
$ws.Cells.Item(34,5).Value = "  -1.91%  "
$c = $ws.Cells.Item(35,4)
$c.NumberFormat = "@"
$c.Value = "2.583"
$c.Style = "Normal"
$ws.Cells.Item(35,5).Value = "  -0.22%  "
$ws.Cells.Item(36,4).Value = "1.232.49"
$ws.Cells.Item(36,5).Value = "  +1.03%  "
$ws.Cells.Item(37,5).Value = "  +0.31%  "
$ws.Cells.Item(38,5).Value = "  -1.80%  "
$c = $ws.Cells.Item(39,4)
$c.NumberFormat = "@"
$c.Value = "6.368"
$c.Style = "Normal"
$ws.Cells.Item(39,5).Value = "  -1.71%  "
$c = $ws.Cells.Item(40,4)
$c.NumberFormat = "@"
$c.Value = "0.9024"
$c.Style = "Normal"
$ws.Cells.Item(40,5).Value = "  -0.61%  "
$c = $ws.Cells.Item(41,4)
$c.NumberFormat = "@"
$c.Value = "0.9986"
$c.Style = "Normal"
$ws.Cells.Item(41,5).Value = "  -0.20%  "
$c = $ws.Cells.Item(42,4)
$c.NumberFormat = "@"
$c.Value = "101.31"
$c.Style = "Normal"
$ws.Cells.Item(42,5).Value = "  -0.31%  "
$c = $ws.Cells.Item(43,4)
$c.NumberFormat = "@"
$c.Value = "65.40"
$c.Style = "Normal"
$ws.Cells.Item(43,5).Value = "  -1.42%  "
$c = $ws.Cells.Item(44,4)
$c.NumberFormat = "@"
$c.Value = "7.164"
$c.Style = "Normal"
$ws.Cells.Item(44,5).Value = "  +0.48%  "
$ws.Cells.Item(45,2).Value = "TheSandbox"
$ws.Cells.Item(45,3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Cells.Item(45,4)
$c.NumberFormat = "@"
$c.Value = "0.3991"
$c.Style = "Normal"
$ws.Cells.Item(45,5).Value = "  -0.70%  "
$ws.Cells.Item(46,2).Value = "EnergySwap"
$ws.Cells.Item(46,3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Cells.Item(46,4)
$c.NumberFormat = "@"
$c.Value = "9.003"
$c.Style = "Normal"
$ws.Cells.Item(46,5).Value = "  -0.11%  "
$ws.Cells.Item(47,2).Value = "RenderToken"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Cells.Item(47,4)
$c.NumberFormat = "@"
$c.Value = "1.684"
$c.Style = "Normal"
$ws.Cells.Item(47,5).Value = "  -0.03%  "
$ws.Cells.Item(48,2).Value = "Algorand"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Cells.Item(48,4)
$c.NumberFormat = "@"
$c.Value = "0.1144"
$c.Style = "Normal"
$ws.Cells.Item(48,5).Value = "  +1.14%  "
$ws.Cells.Item(49,2).Value = "BabyDogeCoin"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c = $ws.Cells.Item(49,4)
$c.NumberFormat = "@"
$c.Value = "0.00000000114"
$c.Style = "Normal"
$ws.Cells.Item(49,5).Value = "  -5.30%  "
$ws.Cells.Item(50,5).Value = "  -0.31%  "
$c = $ws.Cells.Item(51,4)
$c.NumberFormat = "@"
$c.Value = "0.4620"
$c.Style = "Normal"
$ws.Cells.Item(51,5).Value = "  -0.18%  "
